# Update for turkey_super-lig_2023-2024.xlsx
# - Several existing match rows had their home/away (and odds) data swapped
#   between two adjacent rows that share the same match date/time (columns
#   F..V); columns A..E (Indice, pais, torneio, temporada, data_partida) stay
#   put on their own row.
# - Two brand-new match rows (172 and 173) are appended at the end of the
#   sheet, and the sheet's used-range dimension grows from V171 to V173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the F:V ("home" .. "url_partida") contents of two rows, leaving
# A:E (Indice/pais/torneio/temporada/data_partida) untouched on each row.
function Swap-MatchRows($rowA, $rowB) {
    for ($col = 6; $col -le 22; $col++) {
        $valA = $ws.Cells.Item($rowA, $col).Value()
        $valB = $ws.Cells.Item($rowB, $col).Value()
        $ws.Cells.Item($rowA, $col).Value = $valB
        $ws.Cells.Item($rowB, $col).Value = $valA
    }
}

Swap-MatchRows 42 43
Swap-MatchRows 49 50
Swap-MatchRows 53 54
Swap-MatchRows 55 56
Swap-MatchRows 75 76
Swap-MatchRows 127 128
Swap-MatchRows 150 151

# Append two brand-new rows (172 and 173) after the previous last row (171).
# Copy the number formats from row 171's A and E cells (bold/centered/bordered
# index style, and the date-time style) onto the corresponding new cells
# before writing their values, so no extra styles are created.
$ws.Range("A171").Copy()
$ws.Range("A172").PasteSpecial(-4122)
$ws.Range("A173").PasteSpecial(-4122)
$ws.Range("E171").Copy()
$ws.Range("E172").PasteSpecial(-4122)
$ws.Range("E173").PasteSpecial(-4122)

function Set-MatchRow($row, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

Set-MatchRow 172 @(
    171,
    "turkey",
    "super-lig",
    "2023-2024",
    45296.75,
    "Antalyaspor",
    0,
    "Alanyaspor",
    0,
    1.85,
    "28/12/2024 18:43",
    2.12,
    "05/01/2024 17:40",
    3.77,
    "28/12/2024 18:43",
    3.26,
    "05/01/2024 17:40",
    4.25,
    "28/12/2024 18:43",
    4.01,
    "05/01/2024 17:39",
    "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-alanyaspor/IDM0liHn/"
)

Set-MatchRow 173 @(
    172,
    "turkey",
    "super-lig",
    "2023-2024",
    45296.75,
    "Besiktas",
    1,
    "Kasimpasa",
    3,
    1.64,
    "28/12/2024 18:43",
    1.65,
    "05/01/2024 17:57",
    4.28,
    "28/12/2024 18:43",
    4.61,
    "05/01/2024 17:57",
    5,
    "28/12/2024 18:43",
    4.85,
    "05/01/2024 17:59",
    "https://www.betexplorer.com/football/turkey/super-lig/besiktas-kasimpasa/d4L4mBWh/"
)
